$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 553
$ws.Range("J5").Value = 850
$ws.Range("L5").Value = 850
$ws.Range("N5").Value = -1080

$ws.Range("H64").Value = 12362.125
$ws.Range("I64").Value = 7497.5
$ws.Range("K64").Value = 7497.5
$ws.Range("M64").Value = -7249.5

$ws.Range("H67").Value = 12362.125
$ws.Range("I67").Value = 7497.5
$ws.Range("K67").Value = 7497.5
$ws.Range("M67").Value = -6639.5

$ws.Range("H116").Value = 4988.769
$ws.Range("J116").Value = 5420.5
$ws.Range("L116").Value = 5420.5
$ws.Range("N116").Value = -12304.5

$ws.Range("H125").Value = 2583.5
$ws.Range("I125").Value = 2583.5
$ws.Range("K125").Value = 23251.5
$ws.Range("M125").Value = -20791.5


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 954.875
$ws.Range("I2").Value = 973.5
$ws.Range("J2").Value = 948.6667
$ws.Range("K2").Value = 973.5
$ws.Range("L2").Value = 948.6667
$ws.Range("M2").Value = -860.5
$ws.Range("N2").Value = -1174.6667

$ws.Range("H32").Value = 2394.2205
$ws.Range("I32").Value = 943.0682
$ws.Range("J32").Value = 6650.933
$ws.Range("K32").Value = 943.0682
$ws.Range("L32").Value = 6650.933
$ws.Range("M32").Value = -656.0682
$ws.Range("N32").Value = -7224.933

$ws.Range("H45").Value = 64930.25
$ws.Range("I45").Value = 112987.22
$ws.Range("K45").Value = 112987.22
$ws.Range("M45").Value = -112610.22

$ws.Range("H110").Value = 1160.3889
$ws.Range("I110").Value = 1180.8125
$ws.Range("K110").Value = 1180.8125
$ws.Range("M110").Value = 864.1875

$ws.Range("H116").Value = 954.875
$ws.Range("I116").Value = 973.5
$ws.Range("J116").Value = 948.6667
$ws.Range("K116").Value = 973.5
$ws.Range("L116").Value = 948.6667
$ws.Range("M116").Value = 1320.5
$ws.Range("N116").Value = -5536.6667

$ws.Range("H122").Value = 4701.1665
$ws.Range("I122").Value = 4641.4
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 13924.2
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -11474.2
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 2878
$ws.Range("I132").Value = 3013.5386
$ws.Range("K132").Value = 9040.6158
$ws.Range("M132").Value = -6510.6158


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 954.875
$ws.Range("I3").Value = 973.5
$ws.Range("J3").Value = 948.6667
$ws.Range("K3").Value = 973.5
$ws.Range("L3").Value = 948.6667
$ws.Range("M3").Value = -859.5
$ws.Range("N3").Value = -1176.6667

$ws.Range("H20").Value = 789.5909
$ws.Range("I20").Value = 836.5
$ws.Range("K20").Value = 836.5
$ws.Range("M20").Value = -589.5

$ws.Range("H64").Value = 1482.85
$ws.Range("I64").Value = 638
$ws.Range("J64").Value = 1844.9286
$ws.Range("K64").Value = 638
$ws.Range("L64").Value = 1844.9286
$ws.Range("M64").Value = -413
$ws.Range("N64").Value = -2294.9286

$ws.Range("H67").Value = 1482.85
$ws.Range("I67").Value = 638
$ws.Range("J67").Value = 1844.9286
$ws.Range("K67").Value = 638
$ws.Range("L67").Value = 1844.9286
$ws.Range("M67").Value = 142
$ws.Range("N67").Value = -3404.9286

$ws.Range("H86").Value = 15999
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 15999
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 15999
$ws.Range("N86").Value = -18245
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 15999
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 15999
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 79995
$ws.Range("N89").Value = -91227
$ws.Range("M89").ClearContents()

$ws.Range("H134").Value = 26472526
$ws.Range("I134").Value = 1568
$ws.Range("K134").Value = 4704
$ws.Range("M134").Value = -2169


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 28571932
$ws.Range("I16").Value = 28571932
$ws.Range("K16").Value = 28571932
$ws.Range("M16").Value = -28571645

$ws.Range("H33").Value = 4443.6
$ws.Range("I33").Value = 2457.375
$ws.Range("K33").Value = 2457.375
$ws.Range("M33").Value = -2078.375

$ws.Range("H58").Value = 1980.0264
$ws.Range("I58").Value = 2083.6155
$ws.Range("J58").Value = 1926.16
$ws.Range("K58").Value = 2083.6155
$ws.Range("L58").Value = 1926.16
$ws.Range("M58").Value = -1880.6155
$ws.Range("N58").Value = -2332.16

$ws.Range("H94").Value = 1005.15
$ws.Range("I94").Value = 899.6667
$ws.Range("J94").Value = 1023.7647
$ws.Range("K94").Value = 899.6667
$ws.Range("L94").Value = 1023.7647
$ws.Range("M94").Value = -448.6667
$ws.Range("N94").Value = -1925.7647

$ws.Range("H105").Value = 3528.4
$ws.Range("I105").Value = 2580.6667
$ws.Range("K105").Value = 2580.6667
$ws.Range("M105").Value = -833.6667000000002

$ws.Range("H113").Value = 28571932
$ws.Range("I113").Value = 28571932
$ws.Range("K113").Value = 28571932
$ws.Range("M113").Value = -28569762

$ws.Range("H136").Value = 1980.0264
$ws.Range("I136").Value = 2083.6155
$ws.Range("J136").Value = 1926.16
$ws.Range("K136").Value = 6250.8465
$ws.Range("L136").Value = 5778.48
$ws.Range("M136").Value = -3700.8465
$ws.Range("N136").Value = -10878.48


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 480.64285
$ws.Range("J12").Value = 600.2727
$ws.Range("L12").Value = 1800.8181
$ws.Range("N12").Value = -2146.8181


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5403.8335
$ws.Range("I10").Value = 5900
$ws.Range("J10").Value = 5155.75
$ws.Range("K10").Value = 5900
$ws.Range("L10").Value = 5155.75
$ws.Range("M10").Value = -5731
$ws.Range("N10").Value = -5493.75

$ws.Range("I70").Value = 5353.7144
$ws.Range("J70").Value = 5504
$ws.Range("K70").Value = 5353.7144
$ws.Range("L70").Value = 5504
$ws.Range("M70").Value = -5083.7144
$ws.Range("N70").Value = -6044

$ws.Range("I73").Value = 5353.7144
$ws.Range("J73").Value = 5504
$ws.Range("K73").Value = 5353.7144
$ws.Range("L73").Value = 5504
$ws.Range("M73").Value = -4417.7144
$ws.Range("N73").Value = -7376

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 76931
$ws.Range("J109").Value = 88663.75
$ws.Range("L109").Value = 88663.75
$ws.Range("N109").Value = -91437.75


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1850.3
$ws.Range("I132").Value = 1365.931
$ws.Range("K132").Value = 4097.793
$ws.Range("M132").Value = -1567.793

